$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of the document to right
#    after the pPr of the "1 Introduction" heading paragraph (currently
#    paragraph 2, since paragraph 1 is the empty leading paragraph that
#    will be removed below).
# ---------------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

$introRange = $d.Paragraphs(2).Range
$introRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $introRange)

# ---------------------------------------------------------------------
# 2) Remove the empty leading paragraph (Calibri rPr only, no runs).
# ---------------------------------------------------------------------
$d.Paragraphs(1).Range.Delete()

# ---------------------------------------------------------------------
# 3) Merge runs that were needlessly split across identical formatting
#    back into single runs (no visible text change, just de-duplicated
#    <w:r> elements). Each Find/Replace below is scoped with an exact
#    search string that begins and ends precisely on a run boundary so
#    only the intended runs are coalesced.
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "w cost of making bubble teas. Hence opening a new bubble tea shop is a difficult task due to its competiveness. If the shop is opened in the area with high competition, the business will not be feasible. If the shop is opened in the area where ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "w cost of making bubble teas. Hence opening a new bubble tea shop is a difficult task due to its competiveness. If the shop is opened in the area with high competition, the business will not be feasible. If the shop is opened in the area where ",
    2)

$d.Content.Find.Execute(
    "s not so common for the people who lived in that area, the business will not be successful. Therefore, choosing a right location is a huge first step for a successful bubble shop. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "s not so common for the people who lived in that area, the business will not be successful. Therefore, choosing a right location is a huge first step for a successful bubble shop. ",
    2)

$d.Content.Find.Execute(
    "The objective of this capstone project is to find a right location for opening a new bubble tea shop in the city of Toronto, Canada. We will use data from Foursquare API and use Data Science methodologies to analysis the data. We will use clustering in machine learning to come up with a model that will predict that most appropriate location to open a new bubble tea shop.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The objective of this capstone project is to find a right location for opening a new bubble tea shop in the city of Toronto, Canada. We will use data from Foursquare API and use Data Science methodologies to analysis the data. We will use clustering in machine learning to come up with a model that will predict that most appropriate location to open a new bubble tea shop.",
    2)

$d.Content.Find.Execute(
    "Foursquare API to explore various venues from each neighbourhood. There are two ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Foursquare API to explore various venues from each neighbourhood. There are two ",
    2)

$d.Content.Find.Execute(
    "4.2 Examination and Conclusion:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "4.2 Examination and Conclusion:",
    2)

$d.Content.Find.Execute(
    "Cluster 6",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cluster 6",
    2)

Write-Output "done"
